$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128, shifting existing rows 128-154 down to 129-155.
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new record's data.
$ws.Cells.Item(128, 1).Value = 2
$ws.Cells.Item(128, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(128, 3).Value = "Coquimbo"
$ws.Cells.Item(128, 4).Value = 44644
$ws.Cells.Item(128, 5).Value = 4
$ws.Cells.Item(128, 6).Value = 100112043
$ws.Cells.Item(128, 7).Value = "Pepino ensalada"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 500
$ws.Cells.Item(128, 11).Value = 14000
$ws.Cells.Item(128, 12).Value = 16000
$ws.Cells.Item(128, 13).Value = 15000
$ws.Cells.Item(128, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 214
$ws.Cells.Item(128, 17).Value = 70
$ws.Cells.Item(128, 18).Value = "Hortaliza"
